# Add a new data row (row 12) to the active sheet, mirroring the style of
# the existing data rows in column A, then leave the selection where Excel
# would land after typing the new row (cell D12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "tata-15"
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 20

# Column A data cells use the same style (vertical-center + wrap text) as
# the rest of the column; copy it from the cell directly above so the new
# row matches.
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A12").Value = "tata-15"

$ws.Range("D12").Select()
